$d = $word.ActiveDocument

# Step 1: fix typo "ge" -> "je" (creates its own run, separate session)
$rGe = $d.Range(156, 158)
Write-Host "rGe text: [$($rGe.Text)]"
$rGe.Text = "je"
$d.Save()

# Step 2: add trailing space at the end of the paragraph text (separate session -> separate run)
$rEnd = $d.Range(299, 299)
Write-Host "rEnd text: [$($rEnd.Text)]"
$rEnd.InsertBefore(" ")
$d.Save()

Write-Host "Content: [$($d.Content.Text)]"
